$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (170-176) appended to the "Landscaping Data" table.
# Columns: A Date, B Plant_Type, C Plant_Size, D Low, E High, F Temp_Diff(formula),
#          G Rain, H Growth, I Pruned(Quadrant flag "No"), J Quadrant, K Shade,
#          L UV, M Humidity, N Dew_Point, O Pressure, P Wind_Gust, Q Cloud_Cover,
#          R Visibility, S AQI, T Pollen
$rows = @(
    @{R=170; B="Flowering";    C="Large";  D=64; E=83; G=0; H=0.1;  I="No"; J=2; K="Neutral"; L=8; M=0.37; N=54; O=30.1; P=11; Q=0.04; R2=9.9; S=54; T=38},
    @{R=171; B="Nonflowering"; C="Medium"; D=64; E=83; G=0; H=0.2;  I="No"; J=3; K="Bright";  L=8; M=0.37; N=54; O=30.1; P=11; Q=0.04; R2=9.9; S=54; T=38},
    @{R=172; B="Nonflowering"; C="Small";  D=64; E=83; G=0; H=0.1;  I="No"; J=3; K="Bright";  L=8; M=0.37; N=54; O=30.1; P=11; Q=0.04; R2=9.9; S=54; T=38},
    @{R=173; B="Nonflowering"; C="Medium"; D=64; E=83; G=0; H=0.3;  I="No"; J=3; K="Bright";  L=8; M=0.37; N=54; O=30.1; P=11; Q=0.04; R2=9.9; S=54; T=38},
    @{R=174; B="Nonflowering"; C="Medium"; D=64; E=83; G=0; H=0.25; I="No"; J=3; K="Bright";  L=8; M=0.37; N=54; O=30.1; P=11; Q=0.04; R2=9.9; S=54; T=38},
    @{R=175; B="Nonflowering"; C="Large";  D=64; E=83; G=0; H=0.5;  I="No"; J=4; K="Bright";  L=8; M=0.37; N=54; O=30.1; P=11; Q=0.04; R2=9.9; S=54; T=38},
    @{R=176; B="Tree";         C="Medium"; D=64; E=83; G=0; H=1.1;  I="No"; J=1; K="Dark";    L=8; M=0.37; N=54; O=30.1; P=11; Q=0.04; R2=9.9; S=54; T=38}
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Range("F" + $r).Formula = "=ABS(D" + $r + "-E" + $r + ")"
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.R2
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
}

# Copy the date format (m/d/yyyy, the style already used by column A) down
# into the new rows, then (re)write the date value on top of it. Copy +
# PasteSpecial(formats) reuses the existing style index instead of minting a
# new cellXf like a plain NumberFormat assignment would.
$ws.Range("A169").Copy()
$ws.Range("A170:A176").PasteSpecial(-4122)
$excel.CutCopyMode = $false
foreach ($row in $rows) {
    $ws.Cells.Item($row.R, 1).Value = 45811
}

# Scroll the view and update the active selection like the author left it.
$ws.Application.Goto($ws.Range("Q1"), $true)
$ws.Range("U2").Select()
